# Add data for 2022-02-22
# - Rename sheet to reflect new "through" date
# - Update the February row label text
# - Update February and Total row values

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet (sheet name in workbook.xml)
$ws.Name = "Through 2022-02-14"

# Update the February row label (shared string)
$ws.Range("A3").Value = "February (through 02-14)"

# Update February row (row 3) values
$ws.Range("C3").Value = 20
$ws.Range("D3").Value = 34
$ws.Range("E3").Value = 27
$ws.Range("G3").Value = 35
$ws.Range("H3").Value = 68
$ws.Range("I3").Value = 65

# Update Total row (row 4) values
$ws.Range("C4").Value = 71
$ws.Range("D4").Value = 109
$ws.Range("E4").Value = 113
$ws.Range("G4").Value = 109
$ws.Range("H4").Value = 285
$ws.Range("I4").Value = 226
